# Update table from cruise IDs to month names.
# Column C: "OR1-1219" -> "March"
# Column D: "OR1-1242" -> "October"
# Applies to every data row (below the header row) on every worksheet
# that contains these cruise identifiers.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $cCell = $ws.Cells.Item($r, 3)
        $dCell = $ws.Cells.Item($r, 4)

        if ($cCell.Value2 -eq "OR1-1219") {
            $cCell.Value2 = "March"
        }
        if ($dCell.Value2 -eq "OR1-1242") {
            $dCell.Value2 = "October"
        }
    }
}
